# Auto-generated script to apply scheduled-runner price/profit updates
# to the Ravana_Profits workbook. Applies per-cell numeric updates to
# columns H-N (market price / profit calc columns) across all 8 sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3860.25
$ws.Range("J69").Value = 4915
$ws.Range("L69").Value = 14745
$ws.Range("N69").Value = -16493
$ws.Range("H72").Value = 3860.25
$ws.Range("J72").Value = 4915
$ws.Range("L72").Value = 44235
$ws.Range("N72").Value = -52971
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H123").Value = 165390
$ws.Range("J123").Value = 165390
$ws.Range("L123").Value = 165390
$ws.Range("N123").Value = -175190
$ws.Range("H127").Value = 1122.25
$ws.Range("I127").Value = 495
$ws.Range("J127").Value = 1749.5
$ws.Range("K127").Value = 1485
$ws.Range("L127").Value = 5248.5
$ws.Range("M127").Value = 3475
$ws.Range("N127").Value = -15168.5
$ws.Range("H132").Value = 1044.0834
$ws.Range("I132").Value = 1044.0834
$ws.Range("K132").Value = 3132.2502
$ws.Range("M132").Value = -602.2501999999999
$ws.Range("H137").Value = 1855.2667
$ws.Range("I137").Value = 1443.2106
$ws.Range("J137").Value = 2567
$ws.Range("K137").Value = 4329.6318
$ws.Range("L137").Value = 7701
$ws.Range("M137").Value = -1779.6318
$ws.Range("N137").Value = -12801

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H4").Value = 1846.2
$ws.Range("I4").Value = 1829.3334
$ws.Range("K4").Value = 1829.3334
$ws.Range("M4").Value = -1713.3334
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H61").Value = 2657.5386
$ws.Range("I61").Value = 2595.7273
$ws.Range("K61").Value = 2595.7273
$ws.Range("M61").Value = -2383.7273
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 3874.6155
$ws.Range("I132").Value = 2995.8
$ws.Range("K132").Value = 8987.400000000001
$ws.Range("M132").Value = -6457.400000000001
$ws.Range("H136").Value = 2657.5386
$ws.Range("I136").Value = 2595.7273
$ws.Range("K136").Value = 7787.1819
$ws.Range("M136").Value = -5237.1819

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H134").Value = 2617.5454
$ws.Range("I134").Value = 2279.8
$ws.Range("K134").Value = 6839.400000000001
$ws.Range("M134").Value = -4304.400000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H18").Value = 40000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 40000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 40000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -40460
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H107").Value = 2408.2727
$ws.Range("I107").Value = 1750
$ws.Range("J107").Value = 2474.1
$ws.Range("K107").Value = 1750
$ws.Range("L107").Value = 2474.1
$ws.Range("M107").Value = 170
$ws.Range("N107").Value = -6314.1
$ws.Range("H109").Value = 61285
$ws.Range("J109").Value = 61285
$ws.Range("L109").Value = 61285
$ws.Range("N109").Value = -63365
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H114").Value = 47500
$ws.Range("J114").Value = 47500
$ws.Range("L114").Value = 47500
$ws.Range("N114").Value = -56178
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -52350
$ws.Range("H117").Value = 69500
$ws.Range("J117").Value = 69500
$ws.Range("L117").Value = 69500
$ws.Range("N117").Value = -78678
$ws.Range("H120").Value = 29993
$ws.Range("J120").Value = 29993
$ws.Range("L120").Value = 29993
$ws.Range("N120").Value = -37251
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 3351.7778
$ws.Range("I132").Value = 2746.625
$ws.Range("K132").Value = 8239.875
$ws.Range("M132").Value = -5709.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3703932.5
$ws.Range("I2").Value = 4629845.5
$ws.Range("J2").Value = 280
$ws.Range("K2").Value = 27779073
$ws.Range("L2").Value = 1680
$ws.Range("M2").Value = -27778960
$ws.Range("N2").Value = -1906
$ws.Range("H6").Value = 202
$ws.Range("I6").Value = 202
$ws.Range("K6").Value = 606
$ws.Range("M6").Value = -493
$ws.Range("H8").Value = 499.5
$ws.Range("I8").Value = 499.5
$ws.Range("K8").Value = 1498.5
$ws.Range("M8").Value = -1359.5
$ws.Range("H92").Value = 1001.6667
$ws.Range("I92").Value = 1001.6667
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 3005.0001
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1757.0001
$ws.Range("N92").ClearContents()
$ws.Range("H119").Value = 24998.666
$ws.Range("I119").Value = 24998.666
$ws.Range("K119").Value = 74995.99800000001
$ws.Range("M119").Value = -70157.99800000001
$ws.Range("H137").Value = 8128.143
$ws.Range("I137").Value = 10099.75
$ws.Range("K137").Value = 30299.25
$ws.Range("M137").Value = -25199.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4498.6665
$ws.Range("I80").Value = 3397.2
$ws.Range("K80").Value = 3397.2
$ws.Range("M80").Value = -2399.2
$ws.Range("H83").Value = 4498.6665
$ws.Range("I83").Value = 3397.2
$ws.Range("K83").Value = 16986
$ws.Range("M83").Value = -11994
$ws.Range("H97").Value = 498.5625
$ws.Range("I97").Value = 565.2308
$ws.Range("J97").Value = 209.66667
$ws.Range("K97").Value = 565.2308
$ws.Range("L97").Value = 209.66667
$ws.Range("M97").Value = -69.23080000000004
$ws.Range("N97").Value = -1201.66667
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("H113").Value = 1598
$ws.Range("I113").Value = 1598
$ws.Range("K113").Value = 1598
$ws.Range("M113").Value = 572
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3624.5
$ws.Range("I132").Value = 2899.8
$ws.Range("K132").Value = 8699.400000000001
$ws.Range("M132").Value = -6169.400000000001
$ws.Range("H140").Value = 137498.5
$ws.Range("J140").Value = 137498.5
$ws.Range("L140").Value = 137498.5
$ws.Range("N140").Value = -147858.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4500
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 4500
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -4904
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4500
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8840
$ws.Range("H136").Value = 6574.5
$ws.Range("I136").Value = 6574.5
$ws.Range("K136").Value = 19723.5
$ws.Range("M136").Value = -17173.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 10
$ws.Range("K2").Value = 10
$ws.Range("M2").Value = 102
$ws.Range("H96").Value = 5000
$ws.Range("I96").Value = 5000
$ws.Range("K96").Value = 5000
$ws.Range("M96").Value = -3627
$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -480
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 4588.1
$ws.Range("I132").Value = 4172.25
$ws.Range("J132").Value = 4865.3335
$ws.Range("K132").Value = 12516.75
$ws.Range("L132").Value = 14596.0005
$ws.Range("M132").Value = -9986.75
$ws.Range("N132").Value = -19656.0005
$ws.Range("H136").Value = 2678.423
$ws.Range("I136").Value = 2296.8125
$ws.Range("J136").Value = 3289
$ws.Range("K136").Value = 6890.4375
$ws.Range("M136").Value = -4340.4375
$ws.Range("N136").Value = -14967

